# Update TPM-derived values on Sheet1 following recalculation with new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 9.278280333333333
$ws.Range("N2").Value = 27.834841
$ws.Range("O2").Value = 0.2122966588143784
$ws.Range("P2").Value = 0.2122966588143784
$ws.Range("Q2").Value = 0.3687064894062222
$ws.Range("R2").Value = 3.318358404656
$ws.Range("S2").Value = 0.2122966588143784
$ws.Range("T2").Value = 0.2122966588143784

# Row 3
$ws.Range("O3").Value = 0.2154323368929792
$ws.Range("P3").Value = 0.2154323368929792
$ws.Range("S3").Value = 0.2154323368929792
$ws.Range("T3").Value = 0.2154323368929792

# Row 4
$ws.Range("M4").Value = 7.033255
$ws.Range("N4").Value = 21.099765
$ws.Range("O4").Value = 0.1609281551588013
$ws.Range("P4").Value = 0.1609281551588013
$ws.Range("Q4").Value = 0.2794921760266666
$ws.Range("R4").Value = 2.51542958424
$ws.Range("S4").Value = 0.1609281551588013
$ws.Range("T4").Value = 0.1609281551588013

# Row 5
$ws.Range("M5").Value = 17.977458
$ws.Range("N5").Value = 53.932374
$ws.Range("O5").Value = 0.4113428491338411
$ws.Range("P5").Value = 0.411342849133841
$ws.Range("Q5").Value = 0.7144002109759999
$ws.Range("R5").Value = 6.429601898784
$ws.Range("S5").Value = 0.4113428491338411
$ws.Range("T5").Value = 0.411342849133841
